# CRUD for Prospectus Version, Courses, View Student Info 50%
# Adds a YearLvl column (D) and three more sample student rows (3-5),
# including their mailto hyperlinks on the Email column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D
$ws.Range("D1").Value = "YearLvl"

# StudentIDs for the new rows
$ws.Range("A3").Value = "2024-1236"
$ws.Range("A4").Value = "2024-1237"
$ws.Range("A5").Value = "2024-1238"

# Emails (with mailto hyperlinks) for the new rows
$ws.Range("C3").Value = "test1@email.com"
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:test1@email.com")
$ws.Range("C3").Style = "Hyperlink"

$ws.Range("C4").Value = "test2@email.com"
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:test2@email.com")
$ws.Range("C4").Style = "Hyperlink"

$ws.Range("C5").Value = "test3@email.com"
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:test3@email.com")
$ws.Range("C5").Style = "Hyperlink"

# Names for the new rows
$ws.Range("B3").Value = "Test Name 1"
$ws.Range("B4").Value = "Test Name 2"
$ws.Range("B5").Value = "Test Name 3"

# YearLvl values for all data rows
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 4

# Column B widened to fit the new, longer names
$ws.Columns.Item(2).ColumnWidth = 10.33

# Restore active selection to F2, as in the target workbook
$ws.Range("F2").Select() | Out-Null
